$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1197.6923
$ws.Range("I96").Value = 1005.125
$ws.Range("K96").Value = 3015.375
$ws.Range("M96").Value = -1642.375

$ws.Range("H125").Value = 1340.75
$ws.Range("I125").Value = 1436.75
$ws.Range("J125").Value = 1292.75
$ws.Range("K125").Value = 12930.75
$ws.Range("L125").Value = 11634.75
$ws.Range("M125").Value = -10470.75
$ws.Range("N125").Value = -16554.75

$ws.Range("H135").Value = 1900.625
$ws.Range("I135").Value = 1555.9166
$ws.Range("J135").Value = 2934.75
$ws.Range("K135").Value = 14003.2494
$ws.Range("L135").Value = 26412.75
$ws.Range("M135").Value = -11468.2494
$ws.Range("N135").Value = -31482.75

$ws.Range("H138").Value = 8515.244000000001
$ws.Range("J138").Value = 8611.391
$ws.Range("L138").Value = 25834.173
$ws.Range("N138").Value = -36114.173

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22135.666
$ws.Range("I32").Value = 8873.75
$ws.Range("J32").Value = 42538.617
$ws.Range("K32").Value = 8873.75
$ws.Range("L32").Value = 42538.617
$ws.Range("M32").Value = -8586.75
$ws.Range("N32").Value = -43112.617

$ws.Range("H61").Value = 5018.5454
$ws.Range("I61").Value = 4232.316
$ws.Range("K61").Value = 4232.316
$ws.Range("M61").Value = -4020.316

$ws.Range("H63").Value = 4845.591
$ws.Range("I63").Value = 4637.875
$ws.Range("K63").Value = 4637.875
$ws.Range("M63").Value = -3951.875

$ws.Range("H66").Value = 4845.591
$ws.Range("I66").Value = 4637.875
$ws.Range("K66").Value = 23189.375
$ws.Range("M66").Value = -19757.375

$ws.Range("H74").Value = 16744.8
$ws.Range("I74").Value = 22255.715
$ws.Range("K74").Value = 22255.715
$ws.Range("M74").Value = -21381.715

$ws.Range("H77").Value = 16744.8
$ws.Range("I77").Value = 22255.715
$ws.Range("K77").Value = 111278.575
$ws.Range("M77").Value = -106910.575

$ws.Range("H122").Value = 7701.615
$ws.Range("I122").Value = 7849.696
$ws.Range("K122").Value = 23549.088
$ws.Range("M122").Value = -21099.088

$ws.Range("H136").Value = 5018.5454
$ws.Range("I136").Value = 4232.316
$ws.Range("K136").Value = 12696.948
$ws.Range("M136").Value = -10146.948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 16668971
$ws.Range("I99").Value = 16668971
$ws.Range("K99").Value = 16668971
$ws.Range("M99").Value = -16667473

$ws.Range("H135").Value = 73488.5
$ws.Range("J135").Value = 73488.5
$ws.Range("L135").Value = 73488.5
$ws.Range("N135").Value = -83628.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 337510.53
$ws.Range("I31").Value = 3430.6667
$ws.Range("J31").Value = 568796.6
$ws.Range("K31").Value = 3430.6667
$ws.Range("L31").Value = 568796.6
$ws.Range("M31").Value = -3135.6667
$ws.Range("N31").Value = -569386.6

$ws.Range("H34").Value = 337510.53
$ws.Range("I34").Value = 3430.6667
$ws.Range("J34").Value = 568796.6
$ws.Range("K34").Value = 3430.6667
$ws.Range("L34").Value = 568796.6
$ws.Range("M34").Value = -3228.6667
$ws.Range("N34").Value = -569200.6

$ws.Range("H58").Value = 2039.8334
$ws.Range("I58").Value = 1831.238
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 1831.238
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -1628.238
$ws.Range("N58").Value = -3906

$ws.Range("H105").Value = 13287.392
$ws.Range("I105").Value = 13823.182
$ws.Range("K105").Value = 13823.182
$ws.Range("M105").Value = -12076.182

$ws.Range("H136").Value = 2039.8334
$ws.Range("I136").Value = 1831.238
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 5493.714
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2943.714
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 28051.066
$ws.Range("J5").Value = 1609
$ws.Range("L5").Value = 4827
$ws.Range("N5").Value = -5051

$ws.Range("H52").Value = 3985.5
$ws.Range("J52").Value = 3985.5
$ws.Range("L52").Value = 11956.5
$ws.Range("N52").Value = -12488.5

$ws.Range("H113").Value = 20836202
$ws.Range("I113").Value = 3442.7144
$ws.Range("J113").Value = 37039460
$ws.Range("K113").Value = 10328.1432
$ws.Range("L113").Value = 111118380
$ws.Range("M113").Value = -8158.143199999999
$ws.Range("N113").Value = -111122720

$ws.Range("H114").Value = 41313.6
$ws.Range("J114").Value = 102131
$ws.Range("L114").Value = 306393
$ws.Range("N114").Value = -312901

$ws.Range("H129").Value = 7121.3335
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H135").Value = 28051.066
$ws.Range("J135").Value = 1609
$ws.Range("L135").Value = 14481
$ws.Range("N135").Value = -19551

$ws.Range("H137").Value = 8027.5
$ws.Range("J137").Value = 10625
$ws.Range("L137").Value = 31875
$ws.Range("N137").Value = -42075

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 31349.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 31349.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 31349.5
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -32003.5

$ws.Range("H70").Value = 18029.543
$ws.Range("I70").Value = 23270.633
$ws.Range("J70").Value = 8202.5
$ws.Range("K70").Value = 23270.633
$ws.Range("L70").Value = 8202.5
$ws.Range("M70").Value = -23000.633
$ws.Range("N70").Value = -8742.5

$ws.Range("H73").Value = 18029.543
$ws.Range("I73").Value = 23270.633
$ws.Range("J73").Value = 8202.5
$ws.Range("K73").Value = 23270.633
$ws.Range("L73").Value = 8202.5
$ws.Range("M73").Value = -22334.633
$ws.Range("N73").Value = -10074.5

$ws.Range("H80").Value = 5687.375
$ws.Range("I80").Value = 4500
$ws.Range("K80").Value = 4500
$ws.Range("M80").Value = -3502

$ws.Range("H83").Value = 5687.375
$ws.Range("I83").Value = 4500
$ws.Range("K83").Value = 22500
$ws.Range("M83").Value = -17508

$ws.Range("H102").Value = 2981.3438
$ws.Range("I102").Value = 2211.1538
$ws.Range("J102").Value = 6318.8335
$ws.Range("K102").Value = 2211.1538
$ws.Range("L102").Value = 6318.8335
$ws.Range("M102").Value = -589.1538
$ws.Range("N102").Value = -9562.833500000001

$ws.Range("H122").Value = 4184.1
$ws.Range("I122").Value = 4184.1
$ws.Range("K122").Value = 12552.3
$ws.Range("M122").Value = -10102.3

$ws.Range("H126").Value = 5529.074
$ws.Range("I126").Value = 5638.8
$ws.Range("K126").Value = 16916.4
$ws.Range("M126").Value = -14446.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2967.389
$ws.Range("J22").Value = 3065.5293
$ws.Range("L22").Value = 3065.5293
$ws.Range("N22").Value = -3655.5293

$ws.Range("H27").Value = 2967.389
$ws.Range("J27").Value = 3065.5293
$ws.Range("L27").Value = 3065.5293
$ws.Range("N27").Value = -3279.5293

$ws.Range("H68").Value = 4345.4546
$ws.Range("I68").Value = 2387.2222
$ws.Range("K68").Value = 2387.2222
$ws.Range("M68").Value = -1638.2222

$ws.Range("H71").Value = 4345.4546
$ws.Range("I71").Value = 2387.2222
$ws.Range("K71").Value = 11936.111
$ws.Range("M71").Value = -8192.111000000001

$ws.Range("H82").Value = 2428.3
$ws.Range("I82").Value = 1490
$ws.Range("J82").Value = 3835.75
$ws.Range("K82").Value = 1490
$ws.Range("L82").Value = 3835.75
$ws.Range("M82").Value = -1129
$ws.Range("N82").Value = -4557.75

$ws.Range("H85").Value = 2428.3
$ws.Range("I85").Value = 1490
$ws.Range("J85").Value = 3835.75
$ws.Range("K85").Value = 1490
$ws.Range("L85").Value = 3835.75
$ws.Range("M85").Value = -242
$ws.Range("N85").Value = -6331.75

$ws.Range("H136").Value = 9448.416999999999
$ws.Range("I136").Value = 10485.875
$ws.Range("K136").Value = 31457.625
$ws.Range("M136").Value = -28907.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 49832.332
$ws.Range("I54").Value = 49832.332
$ws.Range("K54").Value = 49832.332
$ws.Range("M54").Value = -49312.332

$ws.Range("H107").Value = 1059.4
$ws.Range("I107").Value = 1099.6666
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 3298.9998
$ws.Range("L107").Value = 2997
$ws.Range("M107").Value = -1378.9998
$ws.Range("N107").Value = -6837

$ws.Range("H126").Value = 41669220
$ws.Range("I126").Value = 2480.1
$ws.Range("K126").Value = 7440.299999999999
$ws.Range("M126").Value = -4970.299999999999

$ws.Range("H132").Value = 3249.5334
$ws.Range("I132").Value = 2390.95
$ws.Range("K132").Value = 7172.849999999999
$ws.Range("M132").Value = -4642.849999999999

$ws.Range("H136").Value = 8059.9536
$ws.Range("I136").Value = 6996.423
$ws.Range("K136").Value = 20989.269
$ws.Range("M136").Value = -18439.269
